# Reorders the comma-separated "Recorded By" values in column G so that
# the first entry in the list is moved to the end (rotate left by 1),
# for every row that has more than one entry in that cell.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)
    $val = $cell.Value()

    if ($val -eq $null) {
        continue
    }

    $text = [string]$val
    if ($text -eq "") {
        continue
    }

    $parts = $text.Split(",")
    $trimmed = @()
    foreach ($p in $parts) {
        $trimmed += $p.Trim()
    }

    if ($trimmed.Count -gt 1) {
        $first = $trimmed[0]
        $rest = $trimmed[1..($trimmed.Count - 1)]
        $rotated = $rest + $first
        $joined = [string]::Join(", ", $rotated)
        $cell.Value = $joined
    }
}
